$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue ($ws.Range('D2')) '68.644.56'
$ws.Range('E2').Value = '  +4.26%  '
Set-TextValue ($ws.Range('D3')) '3.370.47'
$ws.Range('E3').Value = '  +1.53%  '
Set-TextValue ($ws.Range('D4')) '1.00'
$ws.Range('E4').Value = '  -0.09%  '
Set-TextValue ($ws.Range('D5')) '594.47'
$ws.Range('E5').Value = '  +6.60%  '
Set-TextValue ($ws.Range('D6')) '186.01'
$ws.Range('E6').Value = '  +0.62%  '
Set-TextValue ($ws.Range('D7')) '0.598'
$ws.Range('E7').Value = '  +4.00%  '
$ws.Range('E8').Value = '  -0.24%  '
Set-TextValue ($ws.Range('D9')) '0.182'
$ws.Range('E9').Value = '  +4.31%  '
Set-TextValue ($ws.Range('D10')) '0.586'
$ws.Range('E10').Value = '  +1.63%  '
Set-TextValue ($ws.Range('D11')) '47.21'
$ws.Range('E11').Value = '  +3.19%  '
Set-TextValue ($ws.Range('D12')) '0.0000280'
$ws.Range('E12').Value = '  +6.98%  '
Set-TextValue ($ws.Range('D13')) '640.57'
$ws.Range('E13').Value = '  +12.72%  '
Set-TextValue ($ws.Range('D14')) '3.915.00'
$ws.Range('E14').Value = '  +1.68%  '
Set-TextValue ($ws.Range('D15')) '8.53'
$ws.Range('E15').Value = '  +1.21%  '
Set-TextValue ($ws.Range('D16')) '68.850.45'
$ws.Range('E16').Value = '  +4.63%  '
$ws.Range('E17').Value = '  +1.93%  '
Set-TextValue ($ws.Range('D18')) '3.376.81'
$ws.Range('E18').Value = '  +1.74%  '
Set-TextValue ($ws.Range('D19')) '17.89'
$ws.Range('E19').Value = '  +1.40%  '
Set-TextValue ($ws.Range('D20')) '11.08'
$ws.Range('E20').Value = '  +2.51%  '
$ws.Range('E21').Value = '  +2.50%  '
Set-TextValue ($ws.Range('D22')) '17.96'
$ws.Range('E22').Value = '  -0.09%  '
Set-TextValue ($ws.Range('D23')) '5.08'
$ws.Range('E23').Value = '  +2.17%  '
Set-TextValue ($ws.Range('D24')) '99.00'
$ws.Range('E24').Value = '  +1.38%  '
Set-TextValue ($ws.Range('D25')) '4.09'
$ws.Range('E25').Value = '  +4.11%  '
Set-TextValue ($ws.Range('D26')) '2.84'
$ws.Range('E26').Value = '  +5.94%  '
$ws.Range('E27').Value = '  +4.98%  '
Set-TextValue ($ws.Range('D28')) '32.92'
$ws.Range('E28').Value = '  +8.29%  '
Set-TextValue ($ws.Range('D29')) '8.67'
$ws.Range('E29').Value = '  +2.81%  '
Set-TextValue ($ws.Range('D30')) '6.80'
$ws.Range('E30').Value = '  +2.07%  '
Set-TextValue ($ws.Range('D31')) '610.12'
$ws.Range('E31').Value = '  +9.08%  '
$ws.Range('B32').Value = 'Maker'
$ws.Range('C32').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue ($ws.Range('D32')) '3.986.22'
$ws.Range('E32').Value = '  +6.75%  '
$ws.Range('B33').Value = 'dogwifhat'
$ws.Range('C33').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue ($ws.Range('D33')) '3.68'
$ws.Range('E33').Value = '  +0.50%  '
Set-TextValue ($ws.Range('D34')) '11.08'
$ws.Range('E34').Value = '  +2.59%  '
$ws.Range('E35').Value = '  +2.68%  '
$ws.Range('E36').Value = '  +0.04%  '
Set-TextValue ($ws.Range('D37')) '56.04'
$ws.Range('E37').Value = '  +0.88%  '
Set-TextValue ($ws.Range('D38')) '2.77'
$ws.Range('E38').Value = '  +7.71%  '
Set-TextValue ($ws.Range('D39')) '3.31'
$ws.Range('E39').Value = '  +6.54%  '
Set-TextValue ($ws.Range('D40')) '0.130'
$ws.Range('E40').Value = '  +4.06%  '
Set-TextValue ($ws.Range('D41')) '33.56'
$ws.Range('E41').Value = '  -0.05%  '
Set-TextValue ($ws.Range('D42')) '0.0₃0705'
$ws.Range('E42').Value = '  +3.37%  '
$ws.Range('E43').Value = '  +2.40%  '
Set-TextValue ($ws.Range('D44')) '0.343'
$ws.Range('E44').Value = '  +3.32%  '
Set-TextValue ($ws.Range('D45')) '0.0422'
$ws.Range('E45').Value = '  +3.55%  '
$ws.Range('E46').Value = '  +2.47%  '
Set-TextValue ($ws.Range('D47')) '2.58'
$ws.Range('E47').Value = '  +3.22%  '
Set-TextValue ($ws.Range('D48')) '1.01'
$ws.Range('E48').Value = '  +0.59%  '
Set-TextValue ($ws.Range('D49')) '1.34'
$ws.Range('E49').Value = '  +9.02%  '
Set-TextValue ($ws.Range('D50')) '131.32'
$ws.Range('E50').Value = '  +4.95%  '
Set-TextValue ($ws.Range('D51')) '7.77'
$ws.Range('E51').Value = '  +7.08%  '
